# Update the "想去人数" (interested-count) figures on the "展览" and
# "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1768
    "F4"  = 5
    "F7"  = 12187
    "F11" = 426
    "F14" = 13565
    "F15" = 13665
    "F23" = 2121
    "F24" = 196
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cell in $updates.Keys) {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
